# Calculated Jan 2020 rent -- fill in Nov/Dec 2019 actuals, add Jan 2020
# actuals and stub out the remaining months of 2020.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$currency = "\$#,##0.00"

# Fill colors (BGR-ordered OLE color values matching the workbook's existing
# themed fills so new cells land on the same cellXfs entries as their
# neighbours instead of minting new styles).
$fillJ = 13553360   # D0CECE-ish fill used by column J (fill index 2)
$fillL = 15983578   # fill used by columns L/S (fill index 3)
$fillT = 14282978   # fill used by column T (fill index 4)
$fillQ = 16777164   # fill used by columns Q/X (fill index 8)

# ---------------------------------------------------------------------
# Row 46 (Sep 2019): the power bill correction.
# ---------------------------------------------------------------------
$ws.Range("D46").Value2 = 63.6

# ---------------------------------------------------------------------
# Row 47 (Oct 2019): record the power bill + Chan's payment.
# ---------------------------------------------------------------------
$ws.Range("D47").NumberFormat = $currency
$ws.Range("D47").Value2 = 45.6

$ws.Range("Q47").NumberFormat = $currency
$ws.Range("Q47").Interior.Color = $fillQ
$ws.Range("Q47").Value2 = 1389.8

# ---------------------------------------------------------------------
# Row 48 (Nov 2019): complete the month's figures.
# ---------------------------------------------------------------------
$ws.Range("D48").NumberFormat = $currency
$ws.Range("D48").Value2 = 147.67

$ws.Range("E48").NumberFormat = $currency
$ws.Range("E48").Formula = "=D47"

$ws.Range("F48").NumberFormat = $currency
$ws.Range("F48").Value2 = 65.99

$ws.Range("G48").NumberFormat = $currency
$ws.Range("G48").Formula = "=C48/3"

$ws.Range("H48").NumberFormat = $currency
$ws.Range("H48").Formula = "=(E48+F48)/3"

$ws.Range("J48").NumberFormat = $currency
$ws.Range("J48").Interior.Color = $fillJ
$ws.Range("J48").Formula = "=G48+H48"

$ws.Range("L48").NumberFormat = $currency
$ws.Range("L48").Interior.Color = $fillL
$ws.Range("L48").Formula = "=3*J48"

$ws.Range("Q48").NumberFormat = $currency
$ws.Range("Q48").Interior.Color = $fillQ
$ws.Range("Q48").Value2 = 424.86

$ws.Range("S48").NumberFormat = $currency
$ws.Range("S48").Interior.Color = $fillL
$ws.Range("S48").Formula = "=S47 + L48 - (J48 + M48+ N48 + O48 + P48 + Q48) + I48"

$ws.Range("T48").NumberFormat = $currency
$ws.Range("T48").Interior.Color = $fillT
$ws.Range("T48").Formula = "=T47 + J48 - M48"

$ws.Range("X48").NumberFormat = $currency
$ws.Range("X48").Interior.Color = $fillQ
$ws.Range("X48").Formula = "=X47 + J48 - Q48 + I48"

# ---------------------------------------------------------------------
# Row 49 (Dec 2019): complete the month's figures (no Chan/Onno payment
# recorded yet).
# ---------------------------------------------------------------------
$ws.Range("D49").NumberFormat = $currency
$ws.Range("D49").Value2 = 149.67

$ws.Range("E49").NumberFormat = $currency
$ws.Range("E49").Formula = "=D48"

$ws.Range("F49").NumberFormat = $currency
$ws.Range("F49").Value2 = 65.99

$ws.Range("G49").NumberFormat = $currency
$ws.Range("G49").Formula = "=C49/3"

$ws.Range("H49").NumberFormat = $currency
$ws.Range("H49").Formula = "=(E49+F49)/3"

$ws.Range("J49").NumberFormat = $currency
$ws.Range("J49").Interior.Color = $fillJ
$ws.Range("J49").Formula = "=G49+H49"

$ws.Range("L49").NumberFormat = $currency
$ws.Range("L49").Interior.Color = $fillL
$ws.Range("L49").Formula = "=3*J49"

$ws.Range("S49").NumberFormat = $currency
$ws.Range("S49").Interior.Color = $fillL
$ws.Range("S49").Formula = "=S48 + L49 - (J49 + M49+ N49 + O49 + P49 + Q49) + I49"

$ws.Range("X49").NumberFormat = $currency
$ws.Range("X49").Interior.Color = $fillQ
$ws.Range("X49").Formula = "=X48 + J49 - Q49 + I49"

# ---------------------------------------------------------------------
# Row 50 (Jan 2020): new year, new month, already calculated.
# ---------------------------------------------------------------------
$ws.Range("A50").Value2 = 2020
$ws.Range("B50").Value = "January"

$ws.Range("C50").NumberFormat = $currency
$ws.Range("C50").Value2 = 1145

$ws.Range("E50").NumberFormat = $currency
$ws.Range("E50").Formula = "=D49"

$ws.Range("F50").NumberFormat = $currency
$ws.Range("F50").Value2 = 65.99

$ws.Range("G50").NumberFormat = $currency
$ws.Range("G50").Formula = "=C50/3"

$ws.Range("H50").NumberFormat = $currency
$ws.Range("H50").Formula = "=(E50+F50)/3"

$ws.Range("J50").NumberFormat = $currency
$ws.Range("J50").Interior.Color = $fillJ
$ws.Range("J50").Formula = "=G50+H50"

$ws.Range("L50").NumberFormat = $currency
$ws.Range("L50").Interior.Color = $fillL
$ws.Range("L50").Formula = "=3*J50"

$ws.Range("S50").NumberFormat = $currency
$ws.Range("S50").Interior.Color = $fillL
$ws.Range("S50").Formula = "=S49 + L50 - (J50 + M50+ N50 + O50 + P50 + Q50) + I50"

$ws.Range("X50").NumberFormat = $currency
$ws.Range("X50").Interior.Color = $fillQ
$ws.Range("X50").Formula = "=X49 + J50 - Q50 + I50"

# ---------------------------------------------------------------------
# Rows 51-61 (Feb-Dec 2020): stub out the rest of the year with the
# known flat rent amount, same as every other future month.
# ---------------------------------------------------------------------
$months = @("February", "March", "April", "May", "June", "July", "August", "September", "October", "November", "December")
$row = 51
foreach ($month in $months) {
    $ws.Range("B$row").Value = $month
    $ws.Range("C$row").NumberFormat = $currency
    $ws.Range("C$row").Value2 = 1145
    $row = $row + 1
}
